# Appends 22 new case rows (752-773) to Sheet1 of the case data workbook,
# matching the data that follows the existing last row (row 751).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$startRow = 752

$data = @(
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00142","Hemmeter","THEFT / M1","2913.02(A)(1)*","M1","Not Guilty"),
    @("22CRB00142","Hemmeter","THEFT / M1","2913.02(A)(1)*","M1","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","DOMESTIC VIOLENCE","2919.25(A)","No Data","Not Guilty"),
    @("22CRB00136","Hemmeter","ASSAULT - M1","2903.13(A)","No Data","Not Guilty")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
